$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Boston_MA
$ws.Range("B2").Value = 0.190191
$ws.Range("C2").Value = 15.499238999999999
$ws.Range("D2").Value = 0
$ws.Range("G2").Value = 0.17596400000000001
$ws.Range("H2").Value = 14.390506
$ws.Range("I2").Value = 0

# Row 3 - Denver_CO
$ws.Range("B3").Value = 0.14471800000000001
$ws.Range("C3").Value = 12.434471
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2.4265729999999999
$ws.Range("F3").Value = 1.2412589999999999
$ws.Range("G3").Value = 0.089472999999999997
$ws.Range("H3").Value = 7.9266959999999997
$ws.Range("I3").Value = 0

# Row 4 - Detroit_MI
$ws.Range("G4").Value = 0.25911099999999998
$ws.Range("H4").Value = 7.3843180000000004
$ws.Range("I4").Value = 0

# Row 5 - Las_Vegas_NV
$ws.Range("G5").Value = 0.62871500000000002
$ws.Range("H5").Value = 18.56551
$ws.Range("I5").Value = 0

# Row 6 - Memphis_TN
$ws.Range("G6").Value = 0.23713699999999999
$ws.Range("H6").Value = 12.060738000000001
$ws.Range("I6").Value = 0

# Update active cell selection from E7 to G7
$ws.Range("G7").Select()
